# Release 1.0.3 renames for the command-line option headers in row 1:
#   -p "Send" -> "Put"
#   -d "Search" -> -s "Search"
# Also the header row got a bit taller to fit the re-wrapped text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Put (-p) String(s)"
$ws.Range("F1").Value = "Search (-s) String(s)"

$ws.Rows("1:1").RowHeight = 51
